$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.935.68"
$ws.Range("E2").Value = "  -6.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.546.74"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.99"
$ws.Range("E6").Value = "  -7.40%  "
$ws.Range("E7").Value = "  -3.93%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -6.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.75"
$ws.Range("E10").Value = "  -8.70%  "
$ws.Range("E11").Value = "  -4.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.68"
$ws.Range("E12").Value = "  -5.86%  "
$ws.Range("E13").Value = "  +0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.938.37"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.544.74"
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.866"
$ws.Range("E16").Value = "  -5.63%  "
$ws.Range("E17").Value = "  -4.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.936.27"
$ws.Range("E18").Value = "  -7.15%  "
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0977"
$ws.Range("E20").Value = "  -4.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.53"
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.05"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "260.70"
$ws.Range("E23").Value = "  -10.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.90"
$ws.Range("E24").Value = "  -5.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "29.46"
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  -7.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("E29").Value = "  -4.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.47"
$ws.Range("E30").Value = "  -7.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.99"
$ws.Range("E31").Value = "  -4.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.67"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("E34").Value = "  -5.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.73"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0792"
$ws.Range("E36").Value = "  -5.57%  "
$ws.Range("E37").Value = "  -6.87%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.119"
$ws.Range("E38").Value = "  -3.64%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.03"
$ws.Range("E39").Value = "  +14.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.75"
$ws.Range("E40").Value = "  +6.30%  "
$ws.Range("E41").Value = "  -3.95%  "
$ws.Range("E42").Value = "  -6.31%  "
$ws.Range("E43").Value = "  -4.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.073.65"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "85.20"
$ws.Range("E46").Value = "  -13.02%  "
$ws.Range("E47").Value = "  +3.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.794.18"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.70"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.16"
$ws.Range("E50").Value = "  -4.35%  "
$ws.Range("E51").Value = "  -9.06%  "
